$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Germanize the title and column headers
$ws.Range("A1").Value = "Bulls Mannschaftsaufstellung: Softball(gemischt) 2013"
$ws.Range("A2").Value = "Vorname"
$ws.Range("B2").Value = "Nachname"
$ws.Range("C2").Value = "Telefon(mobil)"
$ws.Range("D2").Value = "Adresse"
$ws.Range("E2").Value = "Spielposition(en)"

# Column width adjustments
# (target widths: column A -> 13.28515625, column D -> 26.85546875;
#  the inputs below are the values that land closest to those targets
#  once the host engine's internal pixel-grid rounding is applied)
$ws.Columns.Item(1).ColumnWidth = 12.5
$ws.Columns.Item(4).ColumnWidth = 26

# Update the selected cell shown in the sheet view
$ws.Range("D6").Select()
